# [Kadastro App] Yeni kayıt eklendi: 3004
#
# The source workbook stores every column (including the numeric-looking
# "Kayıt No" / "Tarih" / "Parsel Sayısı" columns) as text. Append the new
# record as row 64 to both the master "Kayitlar" sheet and the
# district-filtered "Erdemli" sheet, keeping that same text typing so the
# value isn't silently reinterpreted as a number or a date.

$wb = $excel.ActiveWorkbook

$newRow = @("3004", "2025-09-11", "Erdemli", "1", "ÇAP", "AYHAN KARADAYI (K.Teknisyeni)")
$newRowIndex = 64
$textCols = @(1, 2, 4)   # A: Kayıt No, B: Tarih, D: Parsel Sayısı look numeric/date-like

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($col in $textCols) {
        $ws.Cells.Item($newRowIndex, $col).NumberFormat = "@"
    }

    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item($newRowIndex, $col).Value = $newRow[$col - 1]
    }
}
